# Automatic update of files.
# - Bump the "Förändrad" (C column) date from 45183 to 45184 for every data row.
# - Rewrite the hyperlink formulas in columns S, T, V, W, X, Y for the first
#   four data rows (rows 2-5) to use the two-argument HYPERLINK(url, friendly
#   text) form, and convert the Y column from a plain inline-string
#   (containing an unevaluated, semicolon-separated formula) into a real
#   formula cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump "Förändrad" date for every data row (2 .. 173) -------------
$lastRow = 173
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value = 45184
    }
}

# --- 2. Rewrite hyperlink formulas for rows 2-5 --------------------------
$rowIds = @{
    2 = "A 30234-2023"
    3 = "A 33548-2023"
    4 = "A 33550-2023"
    5 = "A 30241-2023"
}

foreach ($r in ($rowIds.Keys | Sort-Object)) {
    $id = $rowIds[$r]

    # Artfyndslänk (note: source data has a malformed second argument here -
    # missing opening quote before the friendly-name parameter - reproduced
    # verbatim from the upstream export)
    $ws.Range("S$r").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/artfynd/' + $id + '.xlsx, "' + $id + '"")'

    # Kartlänk
    $ws.Range("T$r").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/kartor/' + $id + '.png", "' + $id + '")'

    # Klagomålslänk
    $ws.Range("V$r").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomål/' + $id + '.docx", "' + $id + '")'

    # Klagomålsmaillänk
    $ws.Range("W$r").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/klagomålsmail/' + $id + '.docx", "' + $id + '")'

    # Tillsynsbegäranslänk
    $ws.Range("X$r").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsyn/' + $id + '.docx", "' + $id + '")'

    # Tillsynsbegäransmaillänk (was an inline string with a broken ";"
    # separated formula, now a real formula with a "," separator)
    $ws.Range("Y$r").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ALVDALEN/tillsynsmail/' + $id + '.docx", "' + $id + '")'
}
